$wb = $excel.ActiveWorkbook

# ============================================================
# 1) "总计" summary sheet: insert a new top row for 2022-Q3 and
#    shift the existing quarters down by one.
# ============================================================
$totalWs = $wb.Worksheets.Item("总计")

$totalWs.Rows.Item(2).Insert()

# The inserted row inherited formatting from the row above for B:D;
# strip that back to the plain (unstyled) look used by the other data
# rows, and give the new A2 the same bold/bordered index-column style
# that the rest of column A already uses.
$totalWs.Range("B2:D2").ClearFormats()
$totalWs.Range("A3").Copy()
$totalWs.Range("A2").PasteSpecial(-4122)

$totalWs.Range("A2").Value = 0
$totalWs.Range("B2").Value = "2022-Q3"
$totalWs.Range("C2").Value = 20
$totalWs.Range("D2").Value = 8

# Column A is just a 0-based row index; the rows that got pushed down by
# the insert still carry their old index, so bump each one by 1.
$totalWs.Range("A3").Value = 1
$totalWs.Range("A4").Value = 2
$totalWs.Range("A5").Value = 3
$totalWs.Range("A6").Value = 4
$totalWs.Range("A7").Value = 5
$totalWs.Range("A8").Value = 6

# ============================================================
# 2) New "2022-Q3" sheet with the per-fund holding detail, placed
#    right after "总计" (i.e. right before "2022-Q2").
# ============================================================
$q2Ws = $wb.Worksheets.Item("2022-Q2")
$q2Ws.Copy($q2Ws)
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# The template sheet (2022-Q2) only has 12 data rows; 2022-Q3 needs 20, so
# clone the formatting of the last data row down across the extra 8 rows.
$q3.Range("A13:H13").Copy()
$q3.Range("A14:H21").PasteSpecial(-4122)

# Columns B-G hold text values (fund code / name / percentages as
# strings, matching the source data); force text formatting so the
# numeric-looking ones aren't reinterpreted as numbers.
$q3.Range("B2:G21").NumberFormat = "@"

$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "003567"
$q3.Range("C2").Value = "华夏行业景气混合"
$q3.Range("D2").Value = "115.66"
$q3.Range("E2").Value = "88.33"
$q3.Range("F2").Value = "3.36"
$q3.Range("G2").Value = "3.8862"
$q3.Range("H2").Value = 3
$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "016067"
$q3.Range("C3").Value = "鹏华新能源汽车混合A"
$q3.Range("D3").Value = "21.50"
$q3.Range("E3").Value = "94.32"
$q3.Range("F3").Value = "3.41"
$q3.Range("G3").Value = "0.7332"
$q3.Range("H3").Value = 10
$q3.Range("A4").Value = 2
$q3.Range("B4").Value = "003378"
$q3.Range("C4").Value = "泰康策略优选灵活配置混合"
$q3.Range("D4").Value = "13.78"
$q3.Range("E4").Value = "83.10"
$q3.Range("F4").Value = "4.23"
$q3.Range("G4").Value = "0.5829"
$q3.Range("H4").Value = 5
$q3.Range("A5").Value = 3
$q3.Range("B5").Value = "001349"
$q3.Range("C5").Value = "富国改革动力混合"
$q3.Range("D5").Value = "16.29"
$q3.Range("E5").Value = "77.76"
$q3.Range("F5").Value = "3.38"
$q3.Range("G5").Value = "0.5506"
$q3.Range("H5").Value = 5
$q3.Range("A6").Value = 4
$q3.Range("B6").Value = "010608"
$q3.Range("C6").Value = "华泰柏瑞质量领先混合A"
$q3.Range("D6").Value = "22.70"
$q3.Range("E6").Value = "92.93"
$q3.Range("F6").Value = "1.91"
$q3.Range("G6").Value = "0.4336"
$q3.Range("H6").Value = 4
$q3.Range("A7").Value = 5
$q3.Range("B7").Value = "501202"
$q3.Range("C7").Value = "华泰紫金科技创新3年封闭运作灵活配置混合C"
$q3.Range("D7").Value = "9.12"
$q3.Range("E7").Value = "71.92"
$q3.Range("F7").Value = "4.36"
$q3.Range("G7").Value = "0.3976"
$q3.Range("H7").Value = 3
$q3.Range("A8").Value = 6
$q3.Range("B8").Value = "010874"
$q3.Range("C8").Value = "泰康品质生活混合A"
$q3.Range("D8").Value = "6.86"
$q3.Range("E8").Value = "84.42"
$q3.Range("F8").Value = "4.24"
$q3.Range("G8").Value = "0.2909"
$q3.Range("H8").Value = 7
$q3.Range("A9").Value = 7
$q3.Range("B9").Value = "005825"
$q3.Range("C9").Value = "申万菱信智能驱动股票A"
$q3.Range("D9").Value = "6.11"
$q3.Range("E9").Value = "89.31"
$q3.Range("F9").Value = "4.44"
$q3.Range("G9").Value = "0.2713"
$q3.Range("H9").Value = 4
$q3.Range("A10").Value = 8
$q3.Range("B10").Value = "011769"
$q3.Range("C10").Value = "富国精诚回报12个月持有期混合A"
$q3.Range("D10").Value = "27.25"
$q3.Range("E10").Value = "21.35"
$q3.Range("F10").Value = "0.95"
$q3.Range("G10").Value = "0.2589"
$q3.Range("H10").Value = 7
$q3.Range("A11").Value = 9
$q3.Range("B11").Value = "010875"
$q3.Range("C11").Value = "泰康品质生活混合C"
$q3.Range("D11").Value = "3.44"
$q3.Range("E11").Value = "84.42"
$q3.Range("F11").Value = "4.24"
$q3.Range("G11").Value = "0.1459"
$q3.Range("H11").Value = 7
$q3.Range("A12").Value = 10
$q3.Range("B12").Value = "016068"
$q3.Range("C12").Value = "鹏华新能源汽车混合C"
$q3.Range("D12").Value = "4.13"
$q3.Range("E12").Value = "94.32"
$q3.Range("F12").Value = "3.41"
$q3.Range("G12").Value = "0.1408"
$q3.Range("H12").Value = 10
$q3.Range("A13").Value = 11
$q3.Range("B13").Value = "952035"
$q3.Range("C13").Value = "国泰君安君得诚混合"
$q3.Range("D13").Value = "2.35"
$q3.Range("E13").Value = "80.66"
$q3.Range("F13").Value = "2.92"
$q3.Range("G13").Value = "0.0686"
$q3.Range("H13").Value = 9
$q3.Range("A14").Value = 12
$q3.Range("B14").Value = "010495"
$q3.Range("C14").Value = "创金合信创新驱动股票A"
$q3.Range("D14").Value = "0.95"
$q3.Range("E14").Value = "84.76"
$q3.Range("F14").Value = "7.16"
$q3.Range("G14").Value = "0.0680"
$q3.Range("H14").Value = 1
$q3.Range("A15").Value = 13
$q3.Range("B15").Value = "010609"
$q3.Range("C15").Value = "华泰柏瑞质量领先混合C"
$q3.Range("D15").Value = "2.39"
$q3.Range("E15").Value = "92.93"
$q3.Range("F15").Value = "1.91"
$q3.Range("G15").Value = "0.0456"
$q3.Range("H15").Value = 4
$q3.Range("A16").Value = 14
$q3.Range("B16").Value = "009663"
$q3.Range("C16").Value = "华泰紫金科技创新3年封闭运作灵活配置混合A"
$q3.Range("D16").Value = "0.71"
$q3.Range("E16").Value = "71.92"
$q3.Range("F16").Value = "4.36"
$q3.Range("G16").Value = "0.0310"
$q3.Range("H16").Value = 3
$q3.Range("A17").Value = 15
$q3.Range("B17").Value = "014825"
$q3.Range("C17").Value = "汇泉兴至未来一年持有混合A"
$q3.Range("D17").Value = "1.64"
$q3.Range("E17").Value = "62.93"
$q3.Range("F17").Value = "1.82"
$q3.Range("G17").Value = "0.0298"
$q3.Range("H17").Value = 7
$q3.Range("A18").Value = 16
$q3.Range("B18").Value = "010496"
$q3.Range("C18").Value = "创金合信创新驱动股票C"
$q3.Range("D18").Value = "0.32"
$q3.Range("E18").Value = "84.76"
$q3.Range("F18").Value = "7.16"
$q3.Range("G18").Value = "0.0229"
$q3.Range("H18").Value = 1
$q3.Range("A19").Value = 17
$q3.Range("B19").Value = "015159"
$q3.Range("C19").Value = "申万菱信智能驱动股票C"
$q3.Range("D19").Value = "0.45"
$q3.Range("E19").Value = "89.31"
$q3.Range("F19").Value = "4.44"
$q3.Range("G19").Value = "0.0200"
$q3.Range("H19").Value = 4
$q3.Range("A20").Value = 18
$q3.Range("B20").Value = "011770"
$q3.Range("C20").Value = "富国精诚回报12个月持有期混合C"
$q3.Range("D20").Value = "1.62"
$q3.Range("E20").Value = "21.35"
$q3.Range("F20").Value = "0.95"
$q3.Range("G20").Value = "0.0154"
$q3.Range("H20").Value = 7
$q3.Range("A21").Value = 19
$q3.Range("B21").Value = "014826"
$q3.Range("C21").Value = "汇泉兴至未来一年持有混合C"
$q3.Range("D21").Value = "0.26"
$q3.Range("E21").Value = "62.93"
$q3.Range("F21").Value = "1.82"
$q3.Range("G21").Value = "0.0047"
$q3.Range("H21").Value = 7

